# Apply "Wahl's feedback" changes:
#  1. Swap the contents of columns F and G (data was entered in the wrong columns).
#  2. Clear the File/Folder (A/B) labels on rows that just repeat the label of the
#     row immediately above them (those rows are a continuation of the same file).
#  3. Move the custom column width that belonged to column G onto column F
#     (it travels together with the data that moved from G to F).
#  4. Leave the selection on B6, matching where work resumed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Swap columns F (6) and G (7) for every row that has data ----
$dataRows = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39)

foreach ($r in $dataRows) {
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $fVal = $fCell.Value()
    $gVal = $gCell.Value()

    $fCell.Value = $gVal
    $gCell.Value = $fVal
}

# ---- 2. Clear A/B on rows that merely repeat the row above them ----
$dupRows = @(5,6,8,17,19,21,23,33)

foreach ($r in $dupRows) {
    $ws.Cells.Item($r, 1).ClearContents()
    $ws.Cells.Item($r, 2).ClearContents()
}

# ---- 3. Move the custom width from column G onto column F ----
$gWidth = $ws.Columns.Item(7).ColumnWidth()
$ws.Columns.Item(6).ColumnWidth = $gWidth
$ws.Columns.Item(7).ColumnWidth = 8.43

# ---- 4. Update the active selection ----
$ws.Range("B6").Select()
